$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40 (ALC)
$ws.Range("H40").Value = 999.9286
$ws.Range("I40").Value = 1059.8
$ws.Range("J40").Value = 986.913
$ws.Range("K40").Value = 1059.8
$ws.Range("L40").Value = 986.913
$ws.Range("M40").Value = -884.8
$ws.Range("N40").Value = -1336.913

# Row 51 (ALC)
$ws.Range("H51").Value = 3196.9412
$ws.Range("I51").Value = 1716.6666
$ws.Range("J51").Value = 4004.3635
$ws.Range("K51").Value = 1716.6666
$ws.Range("L51").Value = 4004.3635
$ws.Range("M51").Value = -1232.6666
$ws.Range("N51").Value = -4972.363499999999

# Row 99 (ALC)
$ws.Range("H99").Value = 1021.3
$ws.Range("I99").Value = 401.625
$ws.Range("J99").Value = 3500
$ws.Range("K99").Value = 1204.875
$ws.Range("L99").Value = 10500
$ws.Range("M99").Value = 293.125
$ws.Range("N99").Value = -13496

# Row 113 (ALC)
$ws.Range("H113").Value = 5800
$ws.Range("I113").Value = 5800
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 5800
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2546
$ws.Range("N113").ClearContents()

# Row 138 (ALC)
$ws.Range("H138").Value = 2947665.5
$ws.Range("I138").Value = 128312.125
$ws.Range("J138").Value = 3407968
$ws.Range("K138").Value = 384936.375
$ws.Range("L138").Value = 10223904
$ws.Range("M138").Value = -379796.375
$ws.Range("N138").Value = -10234184

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Range("H2").Value = 1696.2
$ws.Range("I2").Value = 1795.7778
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 1795.7778
$ws.Range("L2").Value = 800
$ws.Range("M2").Value = -1682.7778
$ws.Range("N2").Value = -1026

# Row 4 (ARM)
$ws.Range("H4").Value = 199.5
$ws.Range("I4").Value = 199.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 199.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -83.5
$ws.Range("N4").ClearContents()

# Row 116 (ARM)
$ws.Range("H116").Value = 1696.2
$ws.Range("I116").Value = 1795.7778
$ws.Range("J116").Value = 800
$ws.Range("K116").Value = 1795.7778
$ws.Range("L116").Value = 800
$ws.Range("M116").Value = 498.2221999999999
$ws.Range("N116").Value = -5388

# Row 122 (ARM)
$ws.Range("H122").Value = 3269566.2
$ws.Range("I122").Value = 1640.7778
$ws.Range("J122").Value = 15874422
$ws.Range("K122").Value = 4922.3334
$ws.Range("L122").Value = 47623266
$ws.Range("M122").Value = -2472.3334
$ws.Range("N122").Value = -47628166

# Row 125 (ARM)
$ws.Range("H125").Value = 59799.668
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 59799.668
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 59799.668
$ws.Range("N125").Value = -69639.66800000001

# Row 132 (ARM)
$ws.Range("H132").Value = 97427.336
$ws.Range("I132").Value = 78905.84
$ws.Range("J132").Value = 127524.75
$ws.Range("K132").Value = 236717.52
$ws.Range("L132").Value = 382574.25
$ws.Range("M132").Value = -234187.52
$ws.Range("N132").Value = -387634.25

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Range("H3").Value = 1696.2
$ws.Range("I3").Value = 1795.7778
$ws.Range("J3").Value = 800
$ws.Range("K3").Value = 1795.7778
$ws.Range("L3").Value = 800
$ws.Range("M3").Value = -1681.7778
$ws.Range("N3").Value = -1028

# Row 22 (BSM)
$ws.Range("H22").Value = 680.75
$ws.Range("I22").Value = 637
$ws.Range("J22").Value = 899.5
$ws.Range("K22").Value = 637
$ws.Range("L22").Value = 899.5
$ws.Range("M22").Value = -464
$ws.Range("N22").Value = -1245.5

# Row 35 (BSM)
$ws.Range("H35").Value = 21550
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 21550
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 21550
$ws.Range("M35").Value = -22170

# Row 86 (BSM)
$ws.Range("H86").Value = 11583.625
$ws.Range("I86").Value = 18755.54
$ws.Range("J86").Value = 3107.7273
$ws.Range("K86").Value = 18755.54
$ws.Range("L86").Value = 3107.7273
$ws.Range("M86").Value = -17632.54
$ws.Range("N86").Value = -5353.7273

# Row 89 (BSM)
$ws.Range("H89").Value = 11583.625
$ws.Range("I89").Value = 18755.54
$ws.Range("J89").Value = 3107.7273
$ws.Range("K89").Value = 93777.70000000001
$ws.Range("L89").Value = 15538.6365
$ws.Range("M89").Value = -88161.70000000001
$ws.Range("N89").Value = -26770.6365

# Row 134 (BSM)
$ws.Range("H134").Value = 5484.727
$ws.Range("I134").Value = 4969.6665
$ws.Range("J134").Value = 6102.8
$ws.Range("K134").Value = 14908.9995
$ws.Range("L134").Value = 18308.4
$ws.Range("M134").Value = -12373.9995
$ws.Range("N134").Value = -23378.4

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (CRP)
$ws.Range("H22").Value = 319.15384
$ws.Range("I22").Value = 338.04544
$ws.Range("J22").Value = 215.25
$ws.Range("K22").Value = 338.04544
$ws.Range("L22").Value = 215.25
$ws.Range("M22").Value = 11.95456000000001
$ws.Range("N22").Value = -915.25

# Row 31 (CRP)
$ws.Range("H31").Value = 11787.058
$ws.Range("I31").Value = 37950.535
$ws.Range("J31").Value = 1180.2433
$ws.Range("K31").Value = 37950.535
$ws.Range("L31").Value = 1180.2433
$ws.Range("M31").Value = -37655.535
$ws.Range("N31").Value = -1770.2433

# Row 34 (CRP)
$ws.Range("H34").Value = 11787.058
$ws.Range("I34").Value = 37950.535
$ws.Range("J34").Value = 1180.2433
$ws.Range("K34").Value = 37950.535
$ws.Range("L34").Value = 1180.2433
$ws.Range("M34").Value = -37748.535
$ws.Range("N34").Value = -1584.2433

# Row 62 (CRP)
$ws.Range("H62").Value = 5953.5713
$ws.Range("I62").Value = 4391.6665
$ws.Range("J62").Value = 7125
$ws.Range("K62").Value = 4391.6665
$ws.Range("L62").Value = 7125
$ws.Range("M62").Value = -3767.6665
$ws.Range("N62").Value = -8373

# Row 65 (CRP)
$ws.Range("H65").Value = 5953.5713
$ws.Range("I65").Value = 4391.6665
$ws.Range("J65").Value = 7125
$ws.Range("K65").Value = 21958.3325
$ws.Range("L65").Value = 35625
$ws.Range("M65").Value = -18838.3325
$ws.Range("N65").Value = -41865

# Row 134 (CRP)
$ws.Range("H134").Value = 28314.025
$ws.Range("I134").Value = 1235.7931
$ws.Range("J134").Value = 99702.09
$ws.Range("K134").Value = 3707.379300000001
$ws.Range("L134").Value = 299106.27
$ws.Range("M134").Value = -1172.379300000001
$ws.Range("N134").Value = -304176.27

$ws = $wb.Worksheets.Item("CUL")
# Row 113 (CUL)
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()

# Row 131 (CUL)
$ws.Range("H131").Value = 751.63635
$ws.Range("I131").Value = 559.7143
$ws.Range("J131").Value = 1087.5
$ws.Range("K131").Value = 1679.1429
$ws.Range("L131").Value = 3262.5
$ws.Range("M131").Value = 3360.8571
$ws.Range("N131").Value = -13342.5

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (GSM)
$ws.Range("H80").Value = 3336.2222
$ws.Range("I80").Value = 2447.5
$ws.Range("J80").Value = 4047.2
$ws.Range("K80").Value = 2447.5
$ws.Range("L80").Value = 4047.2
$ws.Range("M80").Value = -1449.5
$ws.Range("N80").Value = -6043.2

# Row 83 (GSM)
$ws.Range("H83").Value = 3336.2222
$ws.Range("I83").Value = 2447.5
$ws.Range("J83").Value = 4047.2
$ws.Range("K83").Value = 12237.5
$ws.Range("L83").Value = 20236
$ws.Range("M83").Value = -7245.5
$ws.Range("N83").Value = -30220

# Row 122 (GSM)
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 3000
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -7900

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (LTW)
$ws.Range("H40").Value = 7628.357
$ws.Range("I40").Value = 7898
$ws.Range("J40").Value = 7143
$ws.Range("K40").Value = 7898
$ws.Range("L40").Value = 7143
$ws.Range("M40").Value = -7762
$ws.Range("N40").Value = -7415

# Row 46 (LTW)
$ws.Range("H46").Value = 775
$ws.Range("I46").Value = 466.66666
$ws.Range("J46").Value = 960
$ws.Range("K46").Value = 466.66666
$ws.Range("L46").Value = 960
$ws.Range("M46").Value = -278.66666
$ws.Range("N46").Value = -1336

# Row 54 (LTW)
$ws.Range("H54").Value = 20000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 20000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 20000
$ws.Range("N54").Value = -21288

# Row 82 (LTW)
$ws.Range("H82").Value = 1819.909
$ws.Range("I82").Value = 1401.8182
$ws.Range("J82").Value = 2238
$ws.Range("K82").Value = 1401.8182
$ws.Range("L82").Value = 2238
$ws.Range("M82").Value = -1040.8182
$ws.Range("N82").Value = -2960

# Row 85 (LTW)
$ws.Range("H85").Value = 1819.909
$ws.Range("I85").Value = 1401.8182
$ws.Range("J85").Value = 2238
$ws.Range("K85").Value = 1401.8182
$ws.Range("L85").Value = 2238
$ws.Range("M85").Value = -153.8181999999999
$ws.Range("N85").Value = -4734

# Row 87 (LTW)
$ws.Range("H87").Value = 38151.2
$ws.Range("I87").Value = 30000
$ws.Range("J87").Value = 40189
$ws.Range("K87").Value = 30000
$ws.Range("L87").Value = 40189
$ws.Range("M87").Value = -28877
$ws.Range("N87").Value = -42435

# Row 90 (LTW)
$ws.Range("H90").Value = 38151.2
$ws.Range("I90").Value = 30000
$ws.Range("J90").Value = 40189
$ws.Range("K90").Value = 90000
$ws.Range("L90").Value = 120567
$ws.Range("M90").Value = -84384
$ws.Range("N90").Value = -131799

# Row 122 (LTW)
$ws.Range("H122").Value = 3938
$ws.Range("I122").Value = 3842.2222
$ws.Range("J122").Value = 4800
$ws.Range("K122").Value = 11526.6666
$ws.Range("L122").Value = 14400
$ws.Range("M122").Value = -9076.6666
$ws.Range("N122").Value = -19300

$ws = $wb.Worksheets.Item("WVR")
# Row 126 (WVR)
$ws.Range("H126").Value = 1547.6333
$ws.Range("I126").Value = 1149.0869
$ws.Range("J126").Value = 2857.1428
$ws.Range("K126").Value = 3447.2607
$ws.Range("L126").Value = 8571.428400000001
$ws.Range("M126").Value = -977.2606999999998
$ws.Range("N126").Value = -13511.4284
